$wb = $excel.ActiveWorkbook

# Sheet "Overview": row 3 is the e460d230 file handback record.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# Sheet "zh-cn": row 3 is the e460d230 file handback record.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-20 02:37:03"

# Sheet "de-de": row 3 is the e460d230 file handback record.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-20 02:37:08"
